# Annot 30 -> 50
# Fill in the "label" (column E) annotation scores for rows 34-55 on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    34 = -1
    35 = 0
    36 = -1
    37 = -2
    38 = 0
    39 = 0
    40 = -1
    41 = 0
    42 = -1
    43 = 0
    44 = 0
    45 = 0
    46 = 0
    47 = 0
    48 = 0
    49 = -1
    50 = -2
    51 = 1
    52 = 0
    53 = 0
    54 = 0
    55 = 0
}

foreach ($row in 34..55) {
    $ws.Cells.Item($row, 5).Value = $values[$row]
}

# Match the author's final selection/scroll position after entering the data.
$ws.Range("E55").Select()
